# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout: 基金代码/基金名称/
#    基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名) to use as the new
#    "2022-Q1" sheet, inserted right before "总计" so it lands in the
#    correct tab order. Duplicating (rather than Worksheets.Add()) keeps
#    the sheetPr/sheetFormatPr/pageMargins/header styling identical to its
#    sibling sheets, then we simply overwrite the two data rows with the
#    2022-Q1 fund figures.
# 2. Rewrite the "总计" (totals) sheet so it gains a new first data row for
#    2022-Q1 (count=2, value=0.12) and the existing rows shift down by one,
#    with the running index in column A renumbered 0..5.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force a numeric-looking string ("010706", "93.50", ...) to be stored
    # as text instead of being auto-coerced to a number (which would drop
    # leading/trailing zeros). Resetting the style back to Normal afterwards
    # avoids leaving a stray "text" number-format style behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------

$template = $wb.Worksheets.Item("2021-Q4")

$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet.Move($null, $template)

# Re-resolve "总计" by name AFTER the sheet insert/move above: worksheet
# handles captured earlier can end up referring to a different sheet once
# the collection is mutated (insertion shifts positions), so look it up
# fresh right before we need it.
$total = $wb.Worksheets.Item("总计")

$fundRows = @(
    @("010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "5.85", "0.0661", 5),
    @("260115", "景顺长城中小盘混合", "0.96", "94.00", "5.22", "0.0501", 6)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. "总计" sheet gains the 2022-Q1 row at the top of the data
# ---------------------------------------------------------------------

$totalsRows = @(
    @("2022-Q1", 2, 0.12),
    @("2021-Q4", 2, 0.42),
    @("2021-Q3", 3, 0.23),
    @("2021-Q2", 2, 0.14),
    @("2021-Q1", 5, 1.21),
    @("2020-Q4", 2, 0.15)
)

# Column A (the running index) carries the same bold/bordered style ("s=2")
# as the header row on every data row. The sheet only had 5 data rows
# before, so row 7 is brand new and needs that style copied onto it
# explicitly (existing rows already carry it).
$indexStyleSource = $total.Cells.Item(2, 1)
$indexStyleSource.Copy()
$total.Cells.Item(7, 1).PasteSpecial(-4122)

for ($i = 0; $i -lt $totalsRows.Length; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]

    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

# Restore the originally-active sheet/selection (sheet insert/copy/move
# operations above shift the active tab as a side effect).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "done"
